$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the specified rows as per repulled/mean-calculated data
$ws.Range("F2").Value = -7
$ws.Range("F3").Value = -3
$ws.Range("F5").Value = -2
$ws.Range("F6").Value = -5
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = -1
$ws.Range("F12").Value = 3
$ws.Range("F14").Value = 2
$ws.Range("F15").Value = -2
